# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets carry identical event data, and the same set of rows need
# their F-column value bumped by the same amount.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1216
    4  = 1127
    5  = 2470
    6  = 8026
    7  = 952
    8  = 494
    10 = 193
    11 = 465
    12 = 21
    13 = 184
    14 = 8402
    16 = 1464
    17 = 171
    20 = 212
    21 = 366
    22 = 213
    23 = 169
    24 = 29
    26 = 40
    27 = 445
    28 = 1193
    29 = 96
    30 = 68
    31 = 115
    32 = 78
    33 = 104
    34 = 54
    35 = 94
    36 = 85
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F" + $row).Value = $updates[$row]
    }
}
